# Update the organization website from "www.stat.kg" to "www.stat.gov.kg"
# on the "Сайт организации (если есть)" row of the metadata sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B10").Select() | Out-Null
